$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Task 9: Register Backend -> status moves from Complete to In Development ---
$ws.Range("G3").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = "In Development"

# --- Task 11: Logout Backend -> status moves from Complete to In Development, add note ---
$ws.Range("G3").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = "In Development"
$ws.Range("C11").Value = "after implementing this, a bug where flash msg persists appeared"

# --- Task 12: Login Backend -> already In Development, add note ---
$ws.Range("C12").Value = "Register, logout, and login has a flash bug, and its error msg needs revamp"

$excel.CutCopyMode = $false

# --- update the active selection / scroll position ---
$ws.Activate()
$ws.Range("C9").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
